$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
if (-not $ws) { $ws = $wb.ActiveSheet }

$ws.Range("B2").Value = 9.174161214183117
$ws.Range("C2").Value = 4.857791381020579
$ws.Range("D2").Value = 4.958651792367739
$ws.Range("F2").Value = 23.82134475437357
$ws.Range("G2").Value = 27.95420364659922
$ws.Range("H2").Value = 14.12660604355934
$ws.Range("I2").Value = 20.03164163482821
$ws.Range("K2").Value = 9.08004166712103
$ws.Range("N2").Value = 17.89298858080832
$ws.Range("B3").Value = 8.863691952697499
$ws.Range("C3").Value = 4.607323665256415
$ws.Range("D3").Value = 4.895570032497602
$ws.Range("F3").Value = 23.86547223680638
$ws.Range("G3").Value = 28.02663593442987
$ws.Range("H3").Value = 14.17754213885025
$ws.Range("I3").Value = 20.12509569139831
$ws.Range("K3").Value = 8.869300602230227
$ws.Range("N3").Value = 17.95124945711763
$ws.Range("B4").Value = 8.668911505669886
$ws.Range("C4").Value = 4.445381584742468
$ws.Range("D4").Value = 4.855908001366819
$ws.Range("F4").Value = 23.89948960618854
$ws.Range("G4").Value = 28.08126294296581
$ws.Range("H4").Value = 14.21126249055108
$ws.Range("I4").Value = 20.1865487851903
$ws.Range("K4").Value = 8.73889279088246
$ws.Range("N4").Value = 17.98876526023162
$ws.Range("B5").Value = 8.588620535939423
$ws.Range("C5").Value = 4.377373768105642
$ws.Range("D5").Value = 4.839520727984005
$ws.Range("F5").Value = 23.91508791782974
$ws.Range("G5").Value = 28.10606428699805
$ws.Range("H5").Value = 14.22561844392268
$ws.Range("I5").Value = 20.21261453066143
$ws.Range("K5").Value = 8.685574494422752
$ws.Range("N5").Value = 18.00449285913779
$ws.Range("B6").Value = 8.575236874607423
$ws.Range("C6").Value = 4.36596058872292
$ws.Range("D6").Value = 4.836786333589527
$ws.Range("F6").Value = 23.91778269666673
$ws.Range("G6").Value = 28.11033556917344
$ws.Range("H6").Value = 14.22803935364686
$ws.Range("I6").Value = 20.21700449832267
$ws.Range("K6").Value = 8.676712747335484
$ws.Range("N6").Value = 18.00713100180323
$ws.Range("B7").Value = 8.667832207423439
$ws.Range("C7").Value = 4.444472514735978
$ws.Range("D7").Value = 4.855687894416655
$ws.Range("F7").Value = 23.89969294862859
$ws.Range("G7").Value = 28.0815871543388
$ws.Range("H7").Value = 14.21145361168473
$ws.Range("I7").Value = 20.18689617545991
$ws.Range("K7").Value = 8.738174328768165
$ws.Range("N7").Value = 17.98897558657664
$ws.Range("B8").Value = 9.068046054378023
$ws.Range("C8").Value = 4.773141445514077
$ws.Range("D8").Value = 4.937099504316357
$ws.Range("F8").Value = 23.83512034917584
$ws.Range("G8").Value = 27.97706369296795
$ws.Range("H8").Value = 14.14366110511468
$ws.Range("I8").Value = 20.06301897482841
$ws.Range("K8").Value = 9.007634268033994
$ws.Range("N8").Value = 17.91271582411034
$ws.Range("B9").Value = 9.814883296741549
$ws.Range("C9").Value = 5.351749181862227
$ws.Range("D9").Value = 5.088986741425787
$ws.Range("F9").Value = 23.76360981423054
$ws.Range("G9").Value = 27.853167414105
$ws.Range("H9").Value = 14.03013547906191
$ws.Range("I9").Value = 19.85244148053885
$ws.Range("K9").Value = 9.524735047721997
$ws.Range("N9").Value = 17.77694753897108
$ws.Range("B10").Value = 10.33444192899218
$ws.Range("C10").Value = 5.735433956751983
$ws.Range("D10").Value = 5.195352519191886
$ws.Range("F10").Value = 23.74487778453306
$ws.Range("G10").Value = 27.81217834842073
$ws.Range("H10").Value = 13.95858095833533
$ws.Range("I10").Value = 19.7174991609852
$ws.Range("K10").Value = 9.893626291666264
$ws.Range("N10").Value = 17.68551777860651
$ws.Range("B11").Value = 10.56341856356021
$ws.Range("C11").Value = 5.900815953714938
$ws.Range("D11").Value = 5.242507354872443
$ws.Range("F11").Value = 23.7437238541116
$ws.Range("G11").Value = 27.80449749384257
$ws.Range("H11").Value = 13.92860641447708
$ws.Range("I11").Value = 19.66041474188778
$ws.Range("K11").Value = 10.0582606963215
$ws.Range("N11").Value = 17.64571355222152
$ws.Range("B12").Value = 10.64899091596742
$ws.Range("C12").Value = 5.962114698639556
$ws.Range("D12").Value = 5.260178651483526
$ws.Range("F12").Value = 23.74434701321476
$ws.Range("G12").Value = 27.80317152355102
$ws.Range("H12").Value = 13.91762662467776
$ws.Range("I12").Value = 19.63941806146146
$ws.Range("K12").Value = 10.1200869036437
$ws.Range("N12").Value = 17.63089654527699
$ws.Range("B13").Value = 10.63061302148037
$ws.Range("C13").Value = 5.948972117108309
$ws.Range("D13").Value = 5.256381184931247
$ws.Range("F13").Value = 23.74416565384137
$ws.Range("G13").Value = 27.80338663332059
$ws.Range("H13").Value = 13.91997481415181
$ws.Range("I13").Value = 19.64391247712578
$ws.Range("K13").Value = 10.10679544499182
$ws.Range("N13").Value = 17.63407628450141
$ws.Range("B14").Value = 10.57048182495041
$ws.Range("C14").Value = 5.90588570265937
$ws.Range("D14").Value = 5.24396494044313
$ws.Range("F14").Value = 23.743753874474
$ws.Range("G14").Value = 27.80435666032764
$ws.Range("H14").Value = 13.92769566391806
$ws.Range("I14").Value = 19.65867489947595
$ws.Range("K14").Value = 10.0633578730652
$ws.Range("N14").Value = 17.64448942380375
$ws.Range("B15").Value = 10.53349967487991
$ws.Range("C15").Value = 5.879320831081724
$ws.Range("D15").Value = 5.236335276096553
$ws.Range("F15").Value = 23.74363971290748
$ws.Range("G15").Value = 27.80515707582575
$ws.Range("H15").Value = 13.9324732268426
$ws.Range("I15").Value = 19.66779808741696
$ws.Range("K15").Value = 10.03668202312536
$ws.Range("N15").Value = 17.65090107821657
$ws.Range("B16").Value = 10.31932270678769
$ws.Range("C16").Value = 5.724440572938248
$ws.Range("D16").Value = 5.192245356382927
$ws.Range("F16").Value = 23.74510152849749
$ws.Range("G16").Value = 27.81290149368891
$ws.Range("H16").Value = 13.96059174284781
$ws.Range("I16").Value = 19.72131643803753
$ws.Range("K16").Value = 9.882797855006412
$ws.Range("N16").Value = 17.68815495837725
$ws.Range("B17").Value = 10.18598673378914
$ws.Range("C17").Value = 5.627072264137655
$ws.Range("D17").Value = 5.16487631327396
$ws.Range("F17").Value = 23.74788606050668
$ws.Range("G17").Value = 27.82046556473497
$ws.Range("H17").Value = 13.97850166218068
$ws.Range("I17").Value = 19.75525091546775
$ws.Range("K17").Value = 9.78753690234111
$ws.Range("N17").Value = 17.71146609696842
$ws.Range("B18").Value = 10.10860534348987
$ws.Range("C18").Value = 5.570207034777974
$ws.Range("D18").Value = 5.149018914912713
$ws.Range("F18").Value = 23.7501811163314
$ws.Range("G18").Value = 27.82584807686194
$ws.Range("H18").Value = 13.98904540617291
$ws.Range("I18").Value = 19.77517392242727
$ws.Range("K18").Value = 9.732449706582212
$ws.Range("N18").Value = 17.72504237057482
$ws.Range("B19").Value = 10.08228927615792
$ws.Range("C19").Value = 5.550805851433652
$ws.Range("D19").Value = 5.143630297128833
$ws.Range("F19").Value = 23.75107724968406
$ws.Range("G19").Value = 27.82784752194884
$ws.Range("H19").Value = 13.99265695954128
$ws.Range("I19").Value = 19.781988992001
$ws.Range("K19").Value = 9.71374927621768
$ws.Range("N19").Value = 17.72966800604318
$ws.Range("B20").Value = 10.20025257319382
$ws.Range("C20").Value = 5.63752654219032
$ws.Range("D20").Value = 5.16780180309442
$ws.Range("F20").Value = 23.7475178630481
$ws.Range("G20").Value = 27.81955351567976
$ws.Range("H20").Value = 13.97657002627307
$ws.Range("I20").Value = 19.75159662785551
$ws.Range("K20").Value = 9.797708641561108
$ws.Range("N20").Value = 17.70896717342836
$ws.Range("B21").Value = 10.5881752094264
$ws.Range("C21").Value = 5.918577337107435
$ws.Range("D21").Value = 5.24761698223561
$ws.Range("F21").Value = 23.74384605142823
$ws.Range("G21").Value = 27.80402875035869
$ws.Range("H21").Value = 13.92541779254991
$ws.Range("I21").Value = 19.65432198157247
$ws.Range("K21").Value = 10.07613103936095
$ws.Range("N21").Value = 17.64142389333834
$ws.Range("B22").Value = 10.83505383334215
$ws.Range("C22").Value = 6.09451747386518
$ws.Range("D22").Value = 5.29869770509876
$ws.Range("F22").Value = 23.74762521804673
$ws.Range("G22").Value = 27.80310856240688
$ws.Range("H22").Value = 13.89414908643427
$ws.Range("I22").Value = 19.59436159909694
$ws.Range("K22").Value = 10.25505625909332
$ws.Range("N22").Value = 17.59877207282979
$ws.Range("B23").Value = 10.70392143539876
$ws.Range("C23").Value = 6.00132628860857
$ws.Range("D23").Value = 5.271536663858258
$ws.Range("F23").Value = 23.74504284040269
$ws.Range("G23").Value = 27.8027539992002
$ws.Range("H23").Value = 13.91063977318484
$ws.Range("I23").Value = 19.62603240967
$ws.Range("K23").Value = 10.15985748140819
$ws.Range("N23").Value = 17.62140001877843
$ws.Range("B24").Value = 10.19380523881283
$ws.Range("C24").Value = 5.632802924733122
$ws.Range("D24").Value = 5.166479570306535
$ws.Range("F24").Value = 23.74768216282967
$ws.Range("G24").Value = 27.81996263322205
$ws.Range("H24").Value = 13.97744254921792
$ws.Range("I24").Value = 19.75324744296861
$ws.Range("K24").Value = 9.79311099159063
$ws.Range("N24").Value = 17.71009639345689
$ws.Range("B25").Value = 9.617569907572873
$ws.Range("C25").Value = 5.202422773511414
$ws.Range("D25").Value = 5.048780229710152
$ws.Range("F25").Value = 23.77703014725986
$ws.Range("G25").Value = 27.8779362565867
$ws.Range("H25").Value = 14.05876692979556
$ws.Range("I25").Value = 19.90594044241988
$ws.Range("K25").Value = 9.386498405155063
$ws.Range("N25").Value = 17.81220982210413
